$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B142").Value = 63902
$ws.Range("E142").Value = 34.04
$ws.Range("F142").Value = 2
$ws.Range("G142").Value = 64.04000000000001

$ws.Range("B143").Value = 48654
$ws.Range("E143").Value = 38.26
$ws.Range("F143").Value = -1
$ws.Range("G143").Value = -32.02

$ws.Range("B154").Value = 64350
$ws.Range("E154").Value = 70.63
$ws.Range("F154").Value = 101
$ws.Range("G154").Value = 6710.44

$ws.Range("B155").Value = 57756
$ws.Range("E155").Value = 79.37
$ws.Range("F155").Value = -100
$ws.Range("G155").Value = -6644

$ws.Range("B156").Value = 53925
$ws.Range("F156").Value = 1
$ws.Range("G156").Value = 66.44

$ws.Range("B256").Value = 64979
$ws.Range("E256").Value = 314.41
$ws.Range("F256").Value = 82
$ws.Range("G256").Value = 24251.5

$ws.Range("B257").Value = 48719
$ws.Range("E257").Value = 353.35
$ws.Range("F257").Value = -81
$ws.Range("G257").Value = -23955.75

$ws.Range("B305").Value = 62997
$ws.Range("F305").Value = 72
$ws.Range("G305").Value = 22020.48

$ws.Range("B306").Value = 57854
$ws.Range("F306").Value = 2
$ws.Range("G306").Value = 611.6799999999999

$ws.Range("B342").Value = 63531
$ws.Range("F342").Value = 80
$ws.Range("G342").Value = 11478.4

$ws.Range("B343").Value = 63571
$ws.Range("E343").Value = 152.53
$ws.Range("F343").Value = 29
$ws.Range("G343").Value = 4160.92

$ws.Range("B344").Value = 57802
$ws.Range("E344").Value = 162.71
$ws.Range("F344").Value = -79
$ws.Range("G344").Value = -11334.92

$ws.Range("B374").Value = 60325
$ws.Range("E374").Value = 151.57
$ws.Range("F374").Value = -102
$ws.Range("G374").Value = -12939.72

$ws.Range("B375").Value = 63560
$ws.Range("E375").Value = 134.87
$ws.Range("F375").Value = 104
$ws.Range("G375").Value = 13193.44

$ws.Range("B381").Value = 62865
$ws.Range("F381").Value = 151
$ws.Range("G381").Value = 12051.31

$ws.Range("B382").Value = 57817
$ws.Range("F382").Value = 3
$ws.Range("G382").Value = 239.43

$ws.Range("B392").Value = 62933
$ws.Range("F392").Value = 146
$ws.Range("G392").Value = 8632.98

$ws.Range("B393").Value = 57835
$ws.Range("F393").Value = 1
$ws.Range("G393").Value = 59.13

$ws.Range("B411").Value = 63007
$ws.Range("F411").Value = 984
$ws.Range("G411").Value = 168588.72

$ws.Range("B412").Value = 57856
$ws.Range("F412").Value = 2
$ws.Range("G412").Value = 342.66

$ws.Range("B413").Value = 63008
$ws.Range("F413").Value = 504
$ws.Range("G413").Value = 76189.67999999999

$ws.Range("B414").Value = 57857
$ws.Range("F414").Value = 3
$ws.Range("G414").Value = 453.51

$ws.Range("B449").Value = 63681
$ws.Range("E449").Value = 23.84
$ws.Range("F449").Value = 65
$ws.Range("G449").Value = 1457.3

$ws.Range("B450").Value = 31930
$ws.Range("E450").Value = 26.8
$ws.Range("F450").Value = -62
$ws.Range("G450").Value = -1390.04

$ws.Range("B578").Value = 64915
$ws.Range("E578").Value = 20.98
$ws.Range("F578").Value = 40
$ws.Range("G578").Value = 789.2

$ws.Range("B579").Value = 45695
$ws.Range("E579").Value = 23.58
$ws.Range("F579").Value = -36
$ws.Range("G579").Value = -710.28

$ws.Range("B596").Value = 53595
$ws.Range("E596").Value = 17.61
$ws.Range("F596").Value = -335
$ws.Range("G596").Value = -4934.55

$ws.Range("B597").Value = 65067
$ws.Range("E597").Value = 15.65
$ws.Range("F597").Value = 338
$ws.Range("G597").Value = 4978.74

$ws.Range("B679").Value = 64810
$ws.Range("E679").Value = 291.22
$ws.Range("F679").Value = 7
$ws.Range("G679").Value = 1917.44

$ws.Range("B680").Value = 53319
$ws.Range("E680").Value = 310.64
$ws.Range("F680").Value = -6
$ws.Range("G680").Value = -1643.52

$ws.Range("B701").Value = 64833
$ws.Range("E701").Value = 34.9
$ws.Range("F701").Value = 99
$ws.Range("G701").Value = 3250.17

$ws.Range("B702").Value = 60025
$ws.Range("E702").Value = 37.22
$ws.Range("F702").Value = -98
$ws.Range("G702").Value = -3217.34

$ws.Range("B712").Value = 64830
$ws.Range("E712").Value = 34.9
$ws.Range("F712").Value = 117
$ws.Range("G712").Value = 3841.11

$ws.Range("B713").Value = 60022
$ws.Range("E713").Value = 37.22
$ws.Range("F713").Value = -113
$ws.Range("G713").Value = -3709.79

$ws.Range("B864").Value = 65079
$ws.Range("E864").Value = 43.44
$ws.Range("F864").Value = 21
$ws.Range("G864").Value = 858.27

$ws.Range("B865").Value = 54751
$ws.Range("E865").Value = 46.34
$ws.Range("F865").Value = -19
$ws.Range("G865").Value = -776.53
